# "feat: add calendar module" - add the next duty day (02.01.2025) to the
# schedule by duplicating the formatting of the last existing row (row 9)
# for "МО \"Большелуцкое сп\"" / Матренина Дарья Анатольевна and bumping
# the date by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 9 (values + styles/borders) into the new row 10 directly below it.
$ws.Range("A9:G9").Copy($ws.Range("A10:G10"))

# The new day is the day after the previous last entry (01.01.2025 -> 02.01.2025).
$ws.Range("B10").Value2 = 45659

# Match the source row's height (rows with a thick bottom border use 15.75).
$ws.Rows.Item(10).RowHeight = 15.75

# Leave the selection where the user would land after typing the new row.
$ws.Range("I11").Select()
